$d = $word.ActiveDocument

$d.Content.Find.Execute("Lynn Willis", $false, $false, $false, $false, $false, $true, 1, $false, "Lynn Willis", 2) | Out-Null
$d.Content.Find.Execute(" in southern Minnesota. His grandfather had emigrated from Norway and founded the town of Frost, Minnesota, where Lon lived after his father died in World War I. His mother later remarried and relocated the family to a farm near Lakefield, Minnesota.", $false, $false, $false, $false, $false, $true, 1, $false, " in southern Minnesota. His grandfather had emigrated from Norway and founded the town of Frost, Minnesota, where Lon lived after his father died in World War I. His mother later remarried and relocated the family to a farm near Lakefield, Minnesota.", 2) | Out-Null
$d.Content.Find.Execute("While in high school, he involved himself in the music and drama departments, and once took second place in the Minnesota State Music Contest playing the saxophone. He won several dramatic competitions, and once won both dramatic and humorous competitions.", $false, $false, $false, $false, $false, $true, 1, $false, "While in high school, he involved himself in the music and drama departments, and once took second place in the Minnesota State Music Contest playing the saxophone. He won several dramatic competitions, and once won both dramatic and humorous competitions.", 2) | Out-Null
$d.Content.Find.Execute("Macthail", $false, $false, $false, $false, $false, $true, 1, $false, "Macthail", 2) | Out-Null
$d.Content.Find.Execute("His first radio appearance was as a soloist in the ", $false, $false, $false, $false, $false, $true, 1, $false, "His first radio appearance was as a soloist in the ", 2) | Out-Null
$d.Content.Find.Execute(" School choir during a Christmas program in 1928. He and a friend teamed up to do musical programs for local radio stations in Minneapolis. Soon after, he joined a tent show that would move to a town, perform different plays for a week, then move on to the next town.", $false, $false, $false, $false, $false, $true, 1, $false, " School choir during a Christmas program in 1928. He and a friend teamed up to do musical programs for local radio stations in Minneapolis. Soon after, he joined a tent show that would move to a town, perform different plays for a week, then move on to the next town.", 2) | Out-Null
$d.Content.Find.Execute(" his talents in radio. He auditioned for station after station with little success. One radio director even told him, “Radio can do very well without you!” Finally, after some time at various Chicago stations, he received an offer from WLW in Cincinnati, Ohio to join a stock company of ten actors. He had the opportunity to play ", $false, $false, $false, $false, $false, $true, 1, $false, " his talents in radio. He auditioned for station after station with little success. One radio director even told him, “Radio can do very well without you!” Finally, after some time at various Chicago stations, he received an offer from WLW in Cincinnati, Ohio to join a stock company of ten actors. He had the opportunity to play ", 2) | Out-Null
$d.Content.Find.Execute(" roles, earning nice reviews for his work.", $false, $false, $false, $false, $false, $true, 1, $false, " roles, earning nice reviews for his work.", 2) | Out-Null
$d.Content.Find.Execute("To further his career, he was faced with a decision: Hollywood with the chance of breaking into films, or New York where he could do legitimate theater. He chose New York, and he and his wife moved there in 1941. He earned close to a thousand dollars in his first month, an unusually large sum for an aspiring actor. Clark was in great demand on the New York radio scene.", $false, $false, $false, $false, $false, $true, 1, $false, "To further his career, he was faced with a decision: Hollywood with the chance of breaking into films, or New York where he could do legitimate theater. He chose New York, and he and his wife moved there in 1941. He earned close to a thousand dollars in his first month, an unusually large sum for an aspiring actor. Clark was in great demand on the New York radio scene.", 2) | Out-Null
$d.Content.Find.Execute("In 1943, he received a call to audition for the Mutual Broadcasting System in a new production based on the Nick Carter dime novels. There were about 50 others auditioning for the role. He was thrilled when Jock MacGregor called him with the news that he had won the leading role for the fledgling series, which lasted until 1955.", $false, $false, $false, $false, $false, $true, 1, $false, "In 1943, he received a call to audition for the Mutual Broadcasting System in a new production based on the Nick Carter dime novels. There were about 50 others auditioning for the role. He was thrilled when Jock MacGregor called him with the news that he had won the leading role for the fledgling series, which lasted until 1955.", 2) | Out-Null
$d.Content.Find.Execute("He also appeared in several other roles in many different shows. During his peak, Lon appeared in an average of twenty radio shows a week, including ", $false, $false, $false, $false, $false, $true, 1, $false, "He also appeared in several other roles in many different shows. During his peak, Lon appeared in an average of twenty radio shows a week, including ", 2) | Out-Null
$d.Content.Find.Execute("Report to the Nation", $false, $false, $false, $false, $false, $true, 1, $false, "Report to the Nation", 2) | Out-Null
$d.Content.Find.Execute(", although he remained uncredited for the role. He did all the voices ", $false, $false, $false, $false, $false, $true, 1, $false, ", although he remained uncredited for the role. He did all the voices ", 2) | Out-Null
$d.Content.Find.Execute(" stood out in Lon Clark’s mind. There was a program that aired during World War II that was sponsored by the Catholic Church, featuring programs about men in the service. He was asked to portray Lt. O’Shay who served in the Navy, and whose ship was sunk with him on board. Prior to the voyage, Lt. O’Shay had written a letter to his son, which was published in newspapers across the country. It was a very emotional Lon Clark who read the letter for the broadcast. The director approached him later, thanking him for saving the show. Ironically, this was the same director who had told the young Lon Clark that “", $false, $false, $false, $false, $false, $true, 1, $false, " stood out in Lon Clark’s mind. There was a program that aired during World War II that was sponsored by the Catholic Church, featuring programs about men in the service. He was asked to portray Lt. O’Shay who served in the Navy, and whose ship was sunk with him on board. Prior to the voyage, Lt. O’Shay had written a letter to his son, which was published in newspapers across the country. It was a very emotional Lon Clark who read the letter for the broadcast. The director approached him later, thanking him for saving the show. Ironically, this was the same director who had told the young Lon Clark that “", 2) | Out-Null
